# Add new rows of data (2023 Serie A round) to BD_Times and BD_Jogo sheets.

$wb = $excel.ActiveWorkbook

$wsTimes = $wb.Worksheets.Item("BD_Times")
$wsJogo  = $wb.Worksheets.Item("BD_Jogo")

# --- BD_Times (sheet1): append rows 416-421, columns A:I ---
$timesData = @(
    @("Athletico PR", 1, 1, 1, 1, 1, 1, 12, 7),
    @("Atletico MG",  0, 1, 1, 1, 1, 1, 7, 12),
    @("Goias",        1, 0, 0, 0, 0, 0, 5, 5),
    @("Internacional",0, 0, 0, 0, 0, 0, 5, 5),
    @("Botafogo",     1, 1, 1, 1, 1, 2, 12, 6),
    @("Flamengo",     0, 1, 1, 1, 2, 1, 6, 12)
)

$startRow = 416
for ($i = 0; $i -lt $timesData.Count; $i++) {
    $r = $startRow + $i
    $row = $timesData[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsTimes.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# --- BD_Jogo (sheet2): append rows 209-211, columns A:E ---
$jogoData = @(
    @(1, 2, 19, "Athletico PR", "Atletico MG"),
    @(0, 0, 10, "Goias", "Internacional"),
    @(1, 3, 18, "Botafogo", "Flamengo")
)

$startRow2 = 209
for ($i = 0; $i -lt $jogoData.Count; $i++) {
    $r = $startRow2 + $i
    $row = $jogoData[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsJogo.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
